# Manual edit on sheet "Card21" (l4.xlsx -> sheet4.xml):
#  - Fill column P ("cylinder(0)") with "nan" placeholder text for every
#    data row (2-34) that was still blank.
#  - Row 9 gets a real service-log entry: Date/Event/Correction/Serviced by
#    are filled in and the new P9 cylinder note is set to "done".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# Row 9: real data entered for this service record.
$ws.Range("L9").Value = "23\1\2026"
$ws.Range("M9").Value = "سيرفيس"
$ws.Range("N9").Value = "تم تغير سلندر"
$ws.Range("O9").Value = "اسطي ايمن تبع مهندس شحته.عيار بواسطه م.عبدالله"
$ws.Range("P9").Value = "done"

# Every other data row (2-34) just gets the "nan" placeholder in column P.
for ($row = 2; $row -le 34; $row++) {
    if ($row -eq 9) {
        continue
    }
    $ws.Cells.Item($row, 16).Value = "nan"
}
